$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the "PLP details" text (shared string reused at the same slot since
# it becomes unreferenced and the new text is appended to the end of the pool)
$ws.Range("B3").Value = "Product List Page Scenarios "

# Turn off wrap-text for the header-ish cell B2 (creates a new cellXfs entry
# that keeps applyAlignment=1 but drops the wrapText override)
$ws.Range("B2").WrapText = $false

# Add the new "Product Details Page" test-case row, re-using the formatting
# from the row above (ProductListPage row) so the style indices line up
$ws.Range("A3:C3").Copy()
$ws.Range("A4:C4").PasteSpecial(-4122)
$ws.Range("B4").Value = "Product Details Page Scenarios"
$ws.Range("A4").Value = "ProductDetailsPage"
$ws.Range("C4").Value = "Y"

# Move the selection to the newly added row, matching where the editor
# left off
$ws.Range("A4").Select() | Out-Null
